# PHPExcel "01simple.xlsx" fixture update (commit 6).
#
# The underlying logical change is: cell A1 on sheet "Simple" now holds
# "                Hello" (16 leading spaces + "Hello") instead of plain
# "Hello". Because the fixture's shared-string table is re-emitted with the
# new string inserted first, every other shared-string index shifts by one
# in the raw XML -- but the actual *values* shown by C1 ("Hello"), B2 and D2
# ("world!") are unchanged. We therefore just (re)write each cell with its
# intended text; the host takes care of shared-string bookkeeping.
#
# The accompanying styles.xml / sheet1.xml churn in the diff (numFmt count,
# font size 10 -> 11, border "outline" flag, cellXfs/cellStyleXfs trimming,
# sheetPr codeName, col width rounding, pageMargins formatting) is PHPExcel
# re-serializing the workbook with a newer writer -- cosmetic/metadata noise
# from a different engine, not a content edit. We still mirror the one
# meaningful, COM-visible piece of it (the base font size) on the cells that
# actually carry data, and make a best-effort (no-op-safe) attempt at the
# sheet codeName, without touching empty cells -- touching whole
# rows/columns/Cells would stamp explicit formatting onto blank cells that
# the target file never touches.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Core content edit -----------------------------------------------
$ws.Range("A1").Value = "                Hello"
$ws.Range("C1").Value = "Hello"
$ws.Range("B2").Value = "world!"
$ws.Range("D2").Value = "world!"

# --- Best-effort formatting touch-ups (kept narrow / no-op safe) ------
# Font size 10 -> 11, applied only to the cells that hold data so blank
# cells in the A1:D2 box stay untouched (matches the target's style 0).
$ws.Range("A1").Font.Size = 11
$ws.Range("C1").Font.Size = 11
$ws.Range("B2").Font.Size = 11
$ws.Range("D2").Font.Size = 11

# sheetPr codeName was dropped by the re-save; CodeName is read-only via
# COM in most hosts, so this is attempted but allowed to silently no-op.
try { $ws.CodeName = "" } catch { }
